$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.379.43'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.839.67'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("E4").Value = '  +1.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.01'
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4746'
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3703'
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07464'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8854'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.52'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.844.12'
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07372'
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.488'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.24'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.581'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008843'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.389.31'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.358'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.081.02'
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.912'
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.33'
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.170'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.255'
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.99'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08967'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7604'
$ws.Range("E32").Value = '  -1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.180'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.558'
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.938'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.013'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05378'
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01960'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.003'
$ws.Range("E40").Value = '  +2.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.288'
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5352'
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.376'
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1665'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.552'
$ws.Range("E45").Value = '  +1.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4974'
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.51'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.014'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.16'
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.679'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06317'
$ws.Range("E51").Value = '  +0.46%  '
